# Update marksheet totals: corrected marks / total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - Right column value
$ws.Range("B11").Value = 5

# "Total" row - Right column value (corrected score)
$ws.Range("B12").Value = 85

# "Total" row - Max column (corrected/total as text string)
$ws.Range("E12").Value = "85/140"
